$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to Price cells that need to preserve exact formatting (trailing zeros, digit grouping)
$priceCells = @(2,3,5,6,8,9,10,11,16,17,19,20,24,30,31,33,39,40,44,46,47,48,50)
foreach ($r in $priceCells) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Update Price (column D) values
$ws.Range("D2").Value = "60.940.86"
$ws.Range("D3").Value = "2.922.30"
$ws.Range("D5").Value = "586.70"
$ws.Range("D6").Value = "145.46"
$ws.Range("D8").Value = "0.506"
$ws.Range("D9").Value = "2.921.63"
$ws.Range("D10").Value = "6.84"
$ws.Range("D11").Value = "0.145"
$ws.Range("D16").Value = "3.406.00"
$ws.Range("D17").Value = "60.906.25"
$ws.Range("D19").Value = "2.919.21"
$ws.Range("D20").Value = "429.17"
$ws.Range("D24").Value = "80.74"
$ws.Range("D30").Value = "7.21"
$ws.Range("D31").Value = "2.62"
$ws.Range("D33").Value = "26.67"
$ws.Range("D39").Value = "0.126"
$ws.Range("D40").Value = "49.65"
$ws.Range("D44").Value = "40.90"
$ws.Range("D46").Value = "379.66"
$ws.Range("D47").Value = "2.693.55"
$ws.Range("D48").Value = "132.55"
$ws.Range("D50").Value = "24.52"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("E3").Value = "  -4.02%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("E6").Value = "  -6.04%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("E9").Value = "  -4.08%  "
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("E13").Value = "  -4.00%  "
$ws.Range("E14").Value = "  -6.54%  "
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("E17").Value = "  -3.54%  "
$ws.Range("E19").Value = "  -4.12%  "
$ws.Range("E20").Value = "  -5.92%  "
$ws.Range("E21").Value = "  -5.04%  "
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  -5.57%  "
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("E26").Value = "  -4.50%  "
$ws.Range("E27").Value = "  -3.26%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  -4.00%  "
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("E37").Value = "  -5.05%  "
$ws.Range("E38").Value = "  -5.67%  "
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("E41").Value = "  -5.76%  "
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("E44").Value = "  -6.67%  "
$ws.Range("E45").Value = "  -2.94%  "
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -2.25%  "
